$d = $word.ActiveDocument

# --- 1. Replace the placeholder body paragraph with the final wording. ---
$d.Content.Find.Execute(
    "Text ", $true, $false, $false, $false, $false, $true, 1, $false,
    "This area of the document is reserved for any requirements not captured in other sections. There are no other requirements at this time.",
    2)

# --- 2. Wrap the (now completed) second paragraph's text in a _Hlk bookmark, ---
#        the artifact Word leaves behind when content is pasted/moved in from
#        elsewhere (this section was moved into the main document).
$bodyPara = $d.Paragraphs.Item(2)
$hlkRange = $d.Range($bodyPara.Range.Start, $bodyPara.Range.End)
$d.Bookmarks.Add("_Hlk480359259", $hlkRange)

# --- 3. Re-seat the _GoBack bookmark (tracks the last edit point) in the ---
#        middle of the heading text, splitting "8. Other requirements" into
#        "8. Other requireme" + "nts" runs around it.
$headingPara = $d.Paragraphs.Item(1)
$splitAt = $headingPara.Range.Start + "8. Other requireme".Length
$goBackRange = $d.Range($splitAt, $splitAt)
$d.Bookmarks.Add("_GoBack", $goBackRange)
